$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("D5", "D6", "D19", "D20", "D24", "D28", "D30", "D32", "D33", "D34", "D35", "D37", "D39", "D42", "D43", "D45", "D46", "D47", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "62.893.25"
$ws.Range("E2").Value = "  -0.96%  "
$ws.Range("D3").Value = "2.544.66"
$ws.Range("E3").Value = "  +3.05%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "567.15"
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("D6").Value = "146.89"
$ws.Range("E6").Value = "  +2.21%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  -1.17%  "
$ws.Range("D9").Value = "2.545.66"
$ws.Range("E9").Value = "  +3.13%  "
$ws.Range("E10").Value = "  -1.43%  "
$ws.Range("E11").Value = "  -2.43%  "
$ws.Range("E12").Value = "  +0.53%  "
$ws.Range("E13").Value = "  -0.53%  "
$ws.Range("D15").Value = "3.002.91"
$ws.Range("E15").Value = "  +3.18%  "
$ws.Range("D16").Value = "62.935.94"
$ws.Range("E16").Value = "  -0.67%  "
$ws.Range("E17").Value = "  -0.29%  "
$ws.Range("D18").Value = "2.540.79"
$ws.Range("E18").Value = "  +2.93%  "
$ws.Range("D19").Value = "11.44"
$ws.Range("E19").Value = "  +1.57%  "
$ws.Range("D20").Value = "335.41"
$ws.Range("E20").Value = "  -1.98%  "
$ws.Range("E21").Value = "  -0.40%  "
$ws.Range("E22").Value = "  -0.92%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "64.92"
$ws.Range("E24").Value = "  -1.23%  "
$ws.Range("E25").Value = "  -3.44%  "
$ws.Range("E26").Value = "  +5.76%  "
$ws.Range("E27").Value = "  +11.43%  "
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.18%  "
$ws.Range("E29").Value = "  +2.82%  "
$ws.Range("D30").Value = "7.26"
$ws.Range("E30").Value = "  +5.40%  "
$ws.Range("D31").Value = "0.0₃0807"
$ws.Range("E31").Value = "  -1.03%  "
$ws.Range("D32").Value = "1.83"
$ws.Range("E32").Value = "  -0.98%  "
$ws.Range("D33").Value = "176.76"
$ws.Range("E33").Value = "  +0.58%  "
$ws.Range("D34").Value = "1.57"
$ws.Range("E34").Value = "  +4.00%  "
$ws.Range("D35").Value = "408.40"
$ws.Range("E35").Value = "  +9.65%  "
$ws.Range("E36").Value = "  -0.63%  "
$ws.Range("D37").Value = "18.96"
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("D39").Value = "4.36"
$ws.Range("E39").Value = "  -2.48%  "
$ws.Range("E40").Value = "  +2.39%  "
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("D42").Value = "39.09"
$ws.Range("E42").Value = "  -3.37%  "
$ws.Range("D43").Value = "153.33"
$ws.Range("E43").Value = "  +1.21%  "
$ws.Range("E44").Value = "  +0.41%  "
$ws.Range("D45").Value = "20.74"
$ws.Range("E45").Value = "  +0.38%  "
$ws.Range("D46").Value = "0.605"
$ws.Range("D47").Value = "0.0956"
$ws.Range("E47").Value = "  -0.71%  "
$ws.Range("E48").Value = "  -0.67%  "
$ws.Range("D49").Value = "0.0235"
$ws.Range("E49").Value = "  +4.09%  "
$ws.Range("D50").Value = "18.21"
$ws.Range("E50").Value = "  +0.78%  "
$ws.Range("E51").Value = "  +0.06%  "
